# "Prix Spot": add a new day column AV (31-jul) with hourly prices for 2025-07-31
$wb = $excel.ActiveWorkbook
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Range("AV1").Value = "31-jul"
$wsPrix.Range("AU1").Copy()
$wsPrix.Range("AV1").PasteSpecial(-4122)

$prixValues = @(90, 77.37, 71.07, 58.38, 57.75, 59.01, 81.96, 100.46, 91.55, 51.31, 30.5, 23.92, 29.01, 23.5, 25.26, 28.5, 30, 42.11, 64.24, 86.71, 107.99, 114.45, 110.27, 98.32)

for ($i = 0; $i -lt $prixValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 48).Value = $prixValues[$i]
}

# "Gaz": add a new row 45 for 2025-07-29
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A45").NumberFormat = "@"
$wsGaz.Range("A45").Value = "2025-07-29"
$wsGaz.Range("A44").Copy()
$wsGaz.Range("A45").PasteSpecial(-4122)
$wsGaz.Range("B45").Value = 33.9

# "CO2": add a new row 45 for 2025-07-29
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A45").NumberFormat = "@"
$wsCo2.Range("A45").Value = "2025-07-29"
$wsCo2.Range("A44").Copy()
$wsCo2.Range("A45").PasteSpecial(-4122)
$wsCo2.Range("B45").Value = 72.16
